$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.318.27'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.596.73'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '1.821.31'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '1.587.24'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.503'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '26.303.43'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.63%  '
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('D20').Value = '0.0₃0720'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.18%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.25%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '1.505.61'
$ws.Range('E32').Value = '  +6.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.567'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.817'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.16'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.935'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.24%  '
$ws.Range('D44').Value = '1.733.90'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.758'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.80%  '
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  +0.07%  '
